# Add a "Feedback" section (heading + paragraph with mailto hyperlink)
# to the end of the document, matching the commit "Add feedback link to pages."

$d = $word.ActiveDocument

# --- Locate the current final paragraph ("Contact details" body paragraph) ---
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

# --- Insert a new paragraph for the "Feedback" heading ---
$insertPoint = $lastPara.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()

$headingParaIndex = $d.Paragraphs.Count
$feedbackHeadingPara = $d.Paragraphs.Item($headingParaIndex)
$feedbackHeadingPara.Range.Text = "Feedback"
$feedbackHeadingPara.Style = "Heading 2"

# --- Insert a new paragraph for the feedback body text ---
$insertPoint2 = $feedbackHeadingPara.Range
$insertPoint2.Collapse(0)
$insertPoint2.InsertParagraphAfter()

$bodyParaIndex = $d.Paragraphs.Count
$feedbackBodyPara = $d.Paragraphs.Item($bodyParaIndex)
$feedbackBodyPara.Style = "Block Text"

# Build the paragraph text in the same run-shape as the source: lead text,
# a space, the e-mail address (will become the hyperlink), then a period.
$r = $feedbackBodyPara.Range
$r.Collapse(0)
$r.InsertAfter("If you have any questions or comments about this guidance, such as suggestions for improvements, please contact:")

$r2 = $feedbackBodyPara.Range
$r2.Collapse(0)
$r2.InsertAfter(" ")

$r3 = $feedbackBodyPara.Range
$r3.Collapse(0)
$r3.InsertAfter("itpolicycontent@digital.justice.gov.uk")

$r4 = $feedbackBodyPara.Range
$r4.Collapse(0)
$r4.InsertAfter(".")

# --- Turn the e-mail address into a mailto: hyperlink ---
$emailRange = $feedbackBodyPara.Range.Duplicate
$emailRange.Find.Execute("itpolicycontent@digital.justice.gov.uk", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($emailRange, "mailto:itpolicycontent@digital.justice.gov.uk", "", "", "itpolicycontent@digital.justice.gov.uk")

# --- Wrap the new "Feedback" heading + paragraph in their own bookmark,
#     mirroring the ariaid-titleN bookmarks used around every other
#     Heading2 section in this document. ---
$wrapStart = $feedbackHeadingPara.Range.Start
$wrapEnd = $feedbackBodyPara.Range.End
$wrapRange = $d.Range($wrapStart, $wrapEnd)
$d.Bookmarks.Add("ariaid-title10", $wrapRange)
